$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.378.94'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '1.837.36'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").Value = "'243.38"
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").Value = "'0.6233"
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D8").Value = "'0.07390"
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("D9").Value = "'0.2930"
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("D10").Value = "'23.24"
$ws.Range("E10").Value = '  +1.04%  '
$ws.Range("D11").Value = "'0.07660"
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").Value = '1.844.70'
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").Value = "'5.012"
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("D14").Value = "'0.6744"
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("D15").Value = "'82.95"
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").Value = "'0.000009198"
$ws.Range("E16").Value = '  +1.62%  '
$ws.Range("D17").Value = "'5.880"
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").Value = '29.378.34'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("D19").Value = '2.096.10'
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").Value = "'238.21"
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").Value = "'12.52"
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("E22").Value = '  -0.56%  '
$ws.Range("D23").Value = "'7.387"
$ws.Range("E23").Value = '  +2.73%  '
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("D25").Value = "'158.42"
$ws.Range("E25").Value = '  -0.87%  '
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("D27").Value = "'8.472"
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("D28").Value = "'17.72"
$ws.Range("E28").Value = '  -0.94%  '
$ws.Range("D29").Value = "'0.06130"
$ws.Range("E29").Value = '  +9.79%  '
$ws.Range("D30").Value = "'1.492"
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("D31").Value = "'1.254"
$ws.Range("E31").Value = '  +3.16%  '
$ws.Range("D32").Value = "'4.115"
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("D33").Value = "'4.092"
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("D34").Value = "'1.857"
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D35").Value = "'1.141"
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("E36").Value = '  -3.33%  '
$ws.Range("D37").Value = "'2.614"
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("D38").Value = "'2.896"
$ws.Range("E38").Value = '  +2.84%  '
$ws.Range("D39").Value = '1.220.17'
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("D40").Value = "'0.01763"
$ws.Range("E40").Value = '  -1.31%  '
$ws.Range("D41").Value = "'6.320"
$ws.Range("E41").Value = '  -1.57%  '
$ws.Range("D42").Value = "'0.9118"
$ws.Range("E42").Value = '  +1.67%  '
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("D44").Value = '2.012.77'
$ws.Range("E44").Value = '  +1.96%  '
$ws.Range("D45").Value = "'101.74"
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("D46").Value = "'65.38"
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("D48").Value = "'0.5072"
$ws.Range("E48").Value = '  -0.81%  '
$ws.Range("D49").Value = "'9.206"
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("D51").Value = "'0.1152"
$ws.Range("E51").Value = '  +3.79%  '
